$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '60.810.39'
$ws.Range('E2').Value = '  -0.23%  '

Set-TextValue $ws.Range('D3') '2.371.25'
$ws.Range('E3').Value = '  -3.87%  '

$ws.Range('E4').Value = '  +0.06%  '

Set-TextValue $ws.Range('D5') '543.62'
$ws.Range('E5').Value = '  -0.74%  '

Set-TextValue $ws.Range('D6') '140.93'
$ws.Range('E6').Value = '  -3.03%  '

$ws.Range('E7').Value = '  +0.03%  '

Set-TextValue $ws.Range('D8') '0.551'
$ws.Range('E8').Value = '  -7.72%  '

Set-TextValue $ws.Range('D9') '2.368.50'
$ws.Range('E9').Value = '  -3.94%  '

$ws.Range('E10').Value = '  -1.85%  '

$ws.Range('E11').Value = '  +0.53%  '

Set-TextValue $ws.Range('D12') '5.33'
$ws.Range('E12').Value = '  -0.50%  '

$ws.Range('E13').Value = '  -2.54%  '

$ws.Range('E14').Value = '  -2.52%  '

Set-TextValue $ws.Range('D15') '2.795.75'
$ws.Range('E15').Value = '  -3.72%  '

$ws.Range('E16').Value = '  +0.03%  '

Set-TextValue $ws.Range('D17') '60.695.40'
$ws.Range('E17').Value = '  -0.24%  '

Set-TextValue $ws.Range('D18') '2.369.15'
$ws.Range('E18').Value = '  -3.50%  '

Set-TextValue $ws.Range('D19') '10.56'
$ws.Range('E19').Value = '  -4.32%  '

$ws.Range('E20').Value = '  -1.91%  '

Set-TextValue $ws.Range('D21') '316.54'
$ws.Range('E21').Value = '  -0.68%  '

$ws.Range('E22').Value = '  -4.00%  '

Set-TextValue $ws.Range('D23') '0.998'
$ws.Range('E23').Value = '  -0.18%  '

Set-TextValue $ws.Range('D24') '1.81'
$ws.Range('E24').Value = '  +3.87%  '

Set-TextValue $ws.Range('D25') '62.83'
$ws.Range('E25').Value = '  -0.54%  '

Set-TextValue $ws.Range('D26') '1.00'
$ws.Range('E26').Value = '  +0.26%  '

Set-TextValue $ws.Range('D27') '2.489.64'
$ws.Range('E27').Value = '  -3.37%  '

$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D28') '7.74'
$ws.Range('E28').Value = '  +1.53%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D29') '0.0₃0923'
$ws.Range('E29').Value = '  -6.11%  '

Set-TextValue $ws.Range('D30') '516.66'
$ws.Range('E30').Value = '  -2.50%  '

$ws.Range('E31').Value = '  -4.16%  '

Set-TextValue $ws.Range('D32') '7.97'
$ws.Range('E32').Value = '  -3.91%  '

$ws.Range('E33').Value = '  -2.99%  '

$ws.Range('E34').Value = '  -3.13%  '

$ws.Range('E35').Value = '  -0.80%  '

Set-TextValue $ws.Range('D36') '0.998'
$ws.Range('E36').Value = '  -0.07%  '

$ws.Range('E37').Value = '  -7.04%  '

$ws.Range('E38').Value = '  -4.64%  '

$ws.Range('E39').Value = '  -0.29%  '

Set-TextValue $ws.Range('D40') '18.03'
$ws.Range('E40').Value = '  -1.49%  '

$ws.Range('E41').Value = '  +1.08%  '

$ws.Range('E42').Value = '  +0.22%  '

Set-TextValue $ws.Range('D43') '136.72'
$ws.Range('E43').Value = '  -5.83%  '

Set-TextValue $ws.Range('D44') '40.21'
$ws.Range('E44').Value = '  +0.82%  '

Set-TextValue $ws.Range('D45') '2.24'
$ws.Range('E45').Value = '  -1.71%  '

Set-TextValue $ws.Range('D46') '138.82'
$ws.Range('E46').Value = '  -5.63%  '

$ws.Range('E47').Value = '  -0.80%  '

Set-TextValue $ws.Range('D48') '20.28'
$ws.Range('E48').Value = '  -2.73%  '

$ws.Range('E49').Value = '  -2.72%  '

Set-TextValue $ws.Range('D50') '0.575'
$ws.Range('E50').Value = '  -1.39%  '

Set-TextValue $ws.Range('D51') '0.0913'
$ws.Range('E51').Value = '  -2.97%  '
